# Rename the worksheet tab (workbook.xml <sheet name="..."> changed to "Sheet1")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Update the data table values (Count column B, Program column D)
$ws.Range("B2").Value = 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "11495"
$ws.Range("B3").Value = 17
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "11495"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "11495"
$ws.Range("B5").Value = 8
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "11495"
$ws.Range("B6").Value = 18
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "11495"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "11495"
$ws.Range("B8").Value = 57
$ws.Range("B9").Value = 1834
$ws.Range("B11").Value = 102
$ws.Range("B12").Value = 1855
$ws.Range("B14").Value = 52
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "143"
$ws.Range("B15").Value = 2621
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "143"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "143"
$ws.Range("B17").Value = 154
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "143"
$ws.Range("B18").Value = 2958
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "143"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "143"
$ws.Range("B20").Value = 167
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8319"
$ws.Range("B21").Value = 5417
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8319"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8319"
$ws.Range("B23").Value = 179
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8319"
$ws.Range("B24").Value = 5668
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8319"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8319"
$ws.Range("B26").Value = 0
$ws.Range("D26").Value = "Erin"
$ws.Range("B27").Value = 75
$ws.Range("D27").Value = "Erin"
$ws.Range("D28").Value = "Erin"
$ws.Range("B29").Value = 1
$ws.Range("D29").Value = "Erin"
$ws.Range("B30").Value = 90
$ws.Range("D30").Value = "Erin"
$ws.Range("D31").Value = "Erin"
$ws.Range("B32").Value = 40
$ws.Range("D32").Value = "MC"
$ws.Range("B33").Value = 697
$ws.Range("D33").Value = "MC"
$ws.Range("D34").Value = "MC"
$ws.Range("B35").Value = 43
$ws.Range("D35").Value = "MC"
$ws.Range("B36").Value = 713
$ws.Range("D36").Value = "MC"
$ws.Range("D37").Value = "MC"
$ws.Range("B38").Value = 41
$ws.Range("D38").Value = "OC"
$ws.Range("B39").Value = 978
$ws.Range("D39").Value = "OC"
$ws.Range("D40").Value = "OC"
$ws.Range("B41").Value = 52
$ws.Range("D41").Value = "OC"
$ws.Range("B42").Value = 995
$ws.Range("D42").Value = "OC"
$ws.Range("D43").Value = "OC"
$ws.Range("B44").Value = 5
$ws.Range("B45").Value = 71
$ws.Range("B47").Value = 22
$ws.Range("B48").Value = 72
